$wb = $excel.ActiveWorkbook

# --- Work on the "pollutant" sheet (sheet1) ---
$wsPollutant = $wb.Worksheets.Item("pollutant")

# Insert a new row 37 so the existing row 37 (PM2.5 -> PM2,5) and everything
# below shifts down by one.
$wsPollutant.Rows.Item(37).Insert()

# Fill in the new row: pollutant_code "NOx" -> pollutant_names "NO2"
$wsPollutant.Range("A37").Value = "NOx"
$wsPollutant.Range("B37").Value = "NO2"

# Grow the pollutant table so the new row is included in it again.
$tbl = $wsPollutant.ListObjects.Item("tbl_pollutant5")
$tbl.Resize($wsPollutant.Range("A1:B46"))

# Select the newly entered cell, matching the author's final selection.
$wsPollutant.Range("B37").Select()

# Make the pollutant sheet the active tab.
$wsPollutant.Activate()
